$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(13).ColumnWidth = 8.166666666666666
$ws.Columns.Item(14).ColumnWidth = 26.619791666666664
